$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark D cells whose new value is a plain number-looking string as Text
# so Excel stores them as literal strings (matching the source data which
# keeps these as text, e.g. "542.44") instead of auto-converting to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.028.46"
$ws.Range("E2").Value = "  +1.52%  "

$ws.Range("D3").Value = "3.095.52"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "542.44"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").Value = "137.36"
$ws.Range("E6").Value = "  +1.38%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "3.089.68"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("D11").Value = "6.43"
$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("E13").Value = "  +6.22%  "

$ws.Range("D14").Value = "34.80"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").Value = "3.595.15"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "64.083.59"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("E17").Value = "  +1.27%  "

$ws.Range("D18").Value = "3.095.20"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "6.71"
$ws.Range("E19").Value = "  +1.72%  "

$ws.Range("D20").Value = "484.53"
$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("D21").Value = "13.41"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  +1.17%  "

$ws.Range("D23").Value = "7.12"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("E24").Value = "  +2.99%  "

$ws.Range("D25").Value = "12.26"
$ws.Range("E25").Value = "  +1.35%  "

$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").Value = "8.12"
$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "26.44"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("E31").Value = "  -0.92%  "

$ws.Range("E32").Value = "  +1.90%  "

$ws.Range("D33").Value = "57.82"
$ws.Range("E33").Value = "  -5.38%  "

$ws.Range("D34").Value = "2.36"
$ws.Range("E34").Value = "  -4.84%  "

$ws.Range("D35").Value = "502.78"
$ws.Range("E35").Value = "  -4.37%  "

$ws.Range("D36").Value = "5.38"
$ws.Range("E36").Value = "  +5.32%  "

$ws.Range("D37").Value = "6.01"
$ws.Range("E37").Value = "  +1.98%  "

$ws.Range("D38").Value = "3.255.41"
$ws.Range("E38").Value = "  +5.68%  "

$ws.Range("D39").Value = "0.0401"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("D40").Value = "0.0797"
$ws.Range("E40").Value = "  +1.81%  "

$ws.Range("E41").Value = "  +1.75%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.70"
$ws.Range("E42").Value = "  +2.33%  "

$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "8.13"
$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("E44").Value = "  +1.60%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.05"
$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "122.99"
$ws.Range("E47").Value = "  +1.55%  "

$ws.Range("D48").Value = "0.0₃0532"
$ws.Range("E48").Value = "  +6.60%  "

$ws.Range("D49").Value = "24.71"
$ws.Range("E49").Value = "  +2.60%  "

$ws.Range("E50").Value = "  +2.43%  "

$ws.Range("D51").Value = "2.41"
$ws.Range("E51").Value = "  +3.53%  "
